# Applies the "Watchdog Timer funktioniert" commit:
#  - Watchdog row (row 21, "Watchdog incl. Vorteiler") gets 2 IST points
#    instead of 0, plus a comment in column F about open questions.
#  - The sheet view scrolls/selects to the newly commented cell (F21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "IST" score for the Watchdog row from 0 to 2.
$ws.Range("D21").Value = 2

# Add the new comment explaining open implementation questions.
$ws.Range("F21").Value = "in Arbeit - die Umsetzung der 18ms sind noch unklar(im Simulator abhängig von Quarzfrequenz?) + ist mit Devise Reset ein Neustart gemeint order ein Stop nach Reset? Und wtf soll ich unter dem Postscaler verstehen?"

# Move the view / selection to the edited cell, scrolled so row 13 is the
# top visible row (mirrors the author's scroll position when saving).
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F21").Select() | Out-Null
